# Split the single "Sheet1" language workbook into five localisation
# sheets (Comm, Property, Guild, Tip, Item) -- matching the authored
# change that combined the per-server language config files into one
# workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet and create the four additional sheets in
#    the right tab order (each new sheet is inserted right after the
#    previous one).
# ---------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item(1)
$wsComm.Name = "Comm"

$wsProperty = $wb.Worksheets.Add($null, $wsComm)
$wsProperty.Name = "Property"

$wsGuild = $wb.Worksheets.Add($null, $wsProperty)
$wsGuild.Name = "Guild"

$wsTip = $wb.Worksheets.Add($null, $wsGuild)
$wsTip.Name = "Tip"

$wsItem = $wb.Worksheets.Add($null, $wsTip)
$wsItem.Name = "Item"

# ---------------------------------------------------------------------
# 2. Comm sheet -- keep the header row, replace the old placeholder rows
#    with the Comm-specific language rows, reusing the existing
#    "Langage_1".."Langage_6" English strings.
# ---------------------------------------------------------------------
$wsComm.Cells.Item(2,1).Value = "Langage_Comm_1"
$wsComm.Cells.Item(2,3).Value = "确认"

$wsComm.Cells.Item(3,1).Value = "Langage_Comm_2"
$wsComm.Cells.Item(3,3).Value = "取消"

$wsComm.Cells.Item(4,1).Value = "Langage_Comm_3"
$wsComm.Cells.Item(4,3).Value = "登录"

$wsComm.Cells.Item(5,1).Value = "Langage_Comm_4"
$wsComm.Cells.Item(5,3).Value = "创建角色"

$wsComm.Cells.Item(6,1).Value = "Langage_Comm_5"
$wsComm.Cells.Item(6,3).Value = "进入游戏"

$wsComm.Cells.Item(7,1).Value = "Langage_Comm_6"
$wsComm.Cells.Item(7,3).Value = "中文_6"

# Five extra blank (format-only) rows below, copied down from row 7.
$wsComm.Range("A7:C7").Copy()
$wsComm.Range("A8:C12").PasteSpecial(-4122)

$wsComm.Columns.Item(1).ColumnWidth = 31.142857142857142
$wsComm.Columns.Item(2).ColumnWidth = 23.714285714285715
$wsComm.Columns.Item(3).ColumnWidth = 22.285714285714285

$wsComm.Range("C8").Select()

# ---------------------------------------------------------------------
# 3. Property sheet -- header row (copied layout) plus six single-column
#    rows of property language keys, then 22 blank format-only rows.
# ---------------------------------------------------------------------
$wsProperty.Cells.Item(1,1).Value = "ID"
$wsProperty.Cells.Item(1,2).Value = "English"
$wsProperty.Cells.Item(1,3).Value = "Chinese"
$wsProperty.Cells.Item(1,2).Font.Bold = $false

$wsProperty.Cells.Item(2,1).Value = "Langage_HP"
$wsProperty.Cells.Item(3,1).Value = "Langage_MAXHP"
$wsProperty.Cells.Item(4,1).Value = "Langage_MP"
$wsProperty.Cells.Item(5,1).Value = "Langage_MAXMP"
$wsProperty.Cells.Item(6,1).Value = "Langage_VP"
$wsProperty.Cells.Item(7,1).Value = "Langage_ATTACK"

# Apply the same formatting as row 1's B/C header cells to the A-only
# data rows (2-7), matching the style used by the rest of the column.
$wsProperty.Range("B1").Copy()
$wsProperty.Range("A2:A7").PasteSpecial(-4122)

# 21 extra blank (format-only) rows below, copied down from row 7.
$wsProperty.Range("A7").Copy()
$wsProperty.Range("A8:A28").PasteSpecial(-4122)

$wsProperty.Columns.Item(1).ColumnWidth = 50.57142857142857

$wsProperty.Range("A1:XFD1").Select()

# ---------------------------------------------------------------------
# 4. Guild sheet -- header row, one data row, then sparse blank rows
#    (3-12, 16) plus a lone formatted cell at A22.
# ---------------------------------------------------------------------
$wsGuild.Cells.Item(1,1).Value = "ID"
$wsGuild.Cells.Item(1,2).Value = "English"
$wsGuild.Cells.Item(1,3).Value = "Chinese"

$wsGuild.Cells.Item(2,1).Value = "Langage_Guild_1"
$wsGuild.Cells.Item(2,2).Value = "Langage_1"
$wsGuild.Cells.Item(2,3).Value = "确认要加入这个公会吗？点击确认加入"

$wsGuild.Range("A2:C2").Copy()
$wsGuild.Range("A3:C12").PasteSpecial(-4122)
$wsGuild.Range("A16:C16").PasteSpecial(-4122)

$wsGuild.Range("B1").Copy()
$wsGuild.Range("A22").PasteSpecial(-4122)

$wsGuild.Columns.Item(1).ColumnWidth = 31.142857142857142
$wsGuild.Columns.Item(2).ColumnWidth = 23.714285714285715
$wsGuild.Columns.Item(3).ColumnWidth = 22.285714285714285

$wsGuild.Range("A12").Select()

# ---------------------------------------------------------------------
# 5. Tip / Item sheets -- header row only.
# ---------------------------------------------------------------------
$wsTip.Cells.Item(1,1).Value = "ID"
$wsTip.Cells.Item(1,2).Value = "English"
$wsTip.Cells.Item(1,3).Value = "Chinese"
$wsTip.Range("A1:XFD1").Select()

$wsItem.Cells.Item(1,1).Value = "ID"
$wsItem.Cells.Item(1,2).Value = "English"
$wsItem.Cells.Item(1,3).Value = "Chinese"
$wsItem.Range("A1:XFD1").Select()

# ---------------------------------------------------------------------
# 6. Leave focus back on the Comm sheet / C8, matching the authored
#    workbook's saved cursor state.
# ---------------------------------------------------------------------
$wsComm.Activate()
$wsComm.Range("C8").Select()

Write-Host "done"
